$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Objects")

# Insert two new rows at position 21 (shifts old rows 21-49 down to 23-51)
$ws.Rows.Item(21).Resize(2).Insert()

$ws.Range("A21").Value = "Image (part)"
